$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 89; this shifts the old rows 89-210 down to 90-211
$ws.Rows(89).Insert()

# Populate the newly inserted row 89 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R keep the same values as the (now shifted) row 90 below it;
# D, J, K, M, P carry the new figures for this record.
$ws.Cells.Item(89, 1).Value = 4
$ws.Cells.Item(89, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(89, 3).Value = "Los Lagos"
$ws.Cells.Item(89, 4).Value = 44579
$ws.Cells.Item(89, 5).Value = 10
$ws.Cells.Item(89, 6).Value = 100112043
$ws.Cells.Item(89, 7).Value = "Pepino ensalada"
$ws.Cells.Item(89, 8).Value = "Sin especificar"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 400
$ws.Cells.Item(89, 11).Value = 14000
$ws.Cells.Item(89, 12).Value = 16000
$ws.Cells.Item(89, 13).Value = 15000
$ws.Cells.Item(89, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(89, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(89, 16).Value = 250
$ws.Cells.Item(89, 17).Value = 60
$ws.Cells.Item(89, 18).Value = "Hortaliza"
